# Update the "想去人数" (interested-count) figures in column F across the
# "展览", "演出" and "全部类型" sheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 264
$ws1.Range("F3").Value = 562
$ws1.Range("F6").Value = 1069
$ws1.Range("F7").Value = 1403
$ws1.Range("F9").Value = 99
$ws1.Range("F12").Value = 129
$ws1.Range("F13").Value = 117
$ws1.Range("F14").Value = 407
$ws1.Range("F15").Value = 1294
$ws1.Range("F16").Value = 97
$ws1.Range("F17").Value = 84
$ws1.Range("F18").Value = 263
$ws1.Range("F20").Value = 634
$ws1.Range("F21").Value = 28
$ws1.Range("F22").Value = 193
$ws1.Range("F23").Value = 8
$ws1.Range("F24").Value = 5565
$ws1.Range("F26").Value = 113
$ws1.Range("F27").Value = 85
$ws1.Range("F29").Value = 14034
$ws1.Range("F30").Value = 1407
$ws1.Range("F31").Value = 189
$ws1.Range("F32").Value = 87
$ws1.Range("F34").Value = 412
$ws1.Range("F35").Value = 565
$ws1.Range("F36").Value = 4156
$ws1.Range("F37").Value = 97

# --- Sheet "演出" (sheet2) ----------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 41

# --- Sheet "全部类型" (sheet4) ------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 264
$ws4.Range("F3").Value = 562
$ws4.Range("F6").Value = 1069
$ws4.Range("F7").Value = 1403
$ws4.Range("F9").Value = 99
$ws4.Range("F12").Value = 129
$ws4.Range("F13").Value = 117
$ws4.Range("F14").Value = 407
$ws4.Range("F15").Value = 1294
$ws4.Range("F16").Value = 97
$ws4.Range("F17").Value = 84
$ws4.Range("F18").Value = 263
$ws4.Range("F21").Value = 634
$ws4.Range("F23").Value = 28
$ws4.Range("F24").Value = 193
$ws4.Range("F25").Value = 8
$ws4.Range("F26").Value = 41
$ws4.Range("F27").Value = 5565
$ws4.Range("F29").Value = 113
$ws4.Range("F30").Value = 85
$ws4.Range("F32").Value = 14035
$ws4.Range("F33").Value = 1407
$ws4.Range("F34").Value = 189
$ws4.Range("F35").Value = 87
$ws4.Range("F37").Value = 412
$ws4.Range("F38").Value = 565
$ws4.Range("F39").Value = 4156
$ws4.Range("F40").Value = 97
